$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$olds = @(
    "77-10=67",
    "32+41=73",
    "94-64=30",
    "48+14=62",
    "47+22=69",
    "81-39=42",
    "99-23=76",
    "45+38=83",
    "47+10=57",
    "30+13=43",
    "1+72=73",
    "74-20=54",
    "42-40=2",
    "53-15=38",
    "85-47=38",
    "80-15=65",
    "73-31=42",
    "47-9=38",
    "50-44=6",
    "76-21=55",
    "58-29=29",
    "92-51=41",
    "49-20=29",
    "92+5=97",
    "31-30=1",
    "89-34=55",
    "21+37=58",
    "8+2=10",
    "95-8=87",
    "2+29=31",
    "26-5=21",
    "65-12=53",
    "61-59=2",
    "82-19=63",
    "37-8=29",
    "23+13=36",
    "91-56=35",
    "86-14=72",
    "92-59=33",
    "84-21=63",
    "9+25=34",
    "39+32=71",
    "92-10=82",
    "39-29=10",
    "48-2=46",
    "74-8=66",
    "41+17=58",
    "56+32=88",
    "44+39=83",
    "77-31=46",
    "43+33=76",
    "96-72=24",
    "59+3=62",
    "36+43=79",
    "16+52=68",
    "97-84=13",
    "58+14=72",
    "79-23=56",
    "84-2=82",
    "12-8=4",
    "28+10=38",
    "83+9=92",
    "13+8=21",
    "16-6=10",
    "90-35=55",
    "29+31=60",
    "4+15=19",
    "69-40=29",
    "16+12=28",
    "48+15=63",
    "72-40=32",
    "42+54=96",
    "25-15=10",
    "47-25=22",
    "37+39=76",
    "90-35=55",
    "3+28=31",
    "77-45=32",
    "76-5=71",
    "6+5=11",
    "32-30=2",
    "77-29=48",
    "57+18=75",
    "74-3=71",
    "87-61=26",
    "72+24=96",
    "93-64=29",
    "55+4=59",
    "48-28=20",
    "65+12=77",
    "47-32=15",
    "35+58=93",
    "26+2=28",
    "74-56=18",
    "60+25=85",
    "21+52=73",
    "37-13=24",
    "18+20=38",
    "83-19=64",
    "19-1=18"
)

$news = @(
    "80-70=10",
    "27+2=29",
    "65-56=9",
    "22+45=67",
    "73+1=74",
    "94-58=36",
    "38+58=96",
    "77+10=87",
    "52-13=39",
    "24-17=7",
    "98-81=17",
    "64+24=88",
    "68-68=0",
    "69+3=72",
    "72+13=85",
    "79-57=22",
    "13+43=56",
    "15+22=37",
    "55+6=61",
    "84-54=30",
    "83-67=16",
    "86-36=50",
    "69-17=52",
    "32-4=28",
    "37+9=46",
    "47+40=87",
    "79+15=94",
    "39-22=17",
    "97-35=62",
    "71+1=72",
    "65-10=55",
    "5+39=44",
    "49+47=96",
    "75-29=46",
    "20+24=44",
    "62-27=35",
    "21+16=37",
    "34-27=7",
    "24+49=73",
    "46-41=5",
    "9+74=83",
    "28+3=31",
    "16-13=3",
    "41-1=40",
    "45+3=48",
    "35+54=89",
    "99-32=67",
    "35+49=84",
    "74-40=34",
    "80+17=97",
    "78-29=49",
    "72-17=55",
    "10+25=35",
    "24+3=27",
    "29+48=77",
    "92-37=55",
    "47-11=36",
    "19+29=48",
    "57+33=90",
    "56+17=73",
    "51-30=21",
    "55-30=25",
    "54-47=7",
    "27+35=62",
    "84-70=14",
    "53-11=42",
    "19-12=7",
    "71-57=14",
    "36+23=59",
    "88-11=77",
    "27+3=30",
    "78+20=98",
    "17+14=31",
    "97-36=61",
    "39+6=45",
    "62-5=57",
    "68-48=20",
    "92-13=79",
    "59-28=31",
    "98-24=74",
    "9-1=8",
    "5+76=81",
    "38-9=29",
    "84+3=87",
    "48+35=83",
    "65+18=83",
    "19-7=12",
    "92-34=58",
    "57+41=98",
    "49+23=72",
    "58+33=91",
    "64-54=10",
    "36+22=58",
    "9+35=44",
    "12+43=55",
    "70+19=89",
    "43-41=2",
    "67-25=42",
    "96-40=56",
    "59-3=56"
)

$cols = 5
$count = $olds.Length
for ($i = 0; $i -lt $count; $i++) {
    $row = [int]([math]::Floor($i / $cols)) + 1
    $col = ($i % $cols) + 1
    $cell = $t.Cell($row, $col)
    $cell.Range.Find.Execute($olds[$i], $true, $false, $false, $false, $false, $true, 1, $false, $news[$i], 1)
}

Write-Host "Replaced" $count "cells"